$wb = $excel.ActiveWorkbook

# --- Step 1: create the new "2022-Q1" sheet by copying "2021-Q4" as a template ---
# (this preserves sheetPr/pageMargins/sheetFormatPr metadata + header style identical to sibling sheets)
$q4 = $wb.Worksheets.Item("2021-Q4")
$totalSheet = $wb.Worksheets.Item("总计")
$q4.Copy($totalSheet)
$q1 = $wb.Worksheets.Item(3)
$q1.Name = "2022-Q1"

# Header row (B1:H1) already correct (copied from 2021-Q4), leave as-is.

# --- Step 2: write the 7 fund rows (rows 2-8) ---

# Row 8 falls beyond the 6 data rows copied from 2021-Q4 (which only had rows 2-7),
# so column A there has no style yet - seed it by copying a sibling index cell first.
$q1.Cells.Item(3, 1).Copy($q1.Cells.Item(8, 1))

# row 2
$q1.Cells.Item(2, 1).Value = 0
$q1.Cells.Item(2, 2).Value = "'005825"
$q1.Cells.Item(2, 2).Style = "Normal"
$q1.Cells.Item(2, 3).Value = "'申万菱信智能驱动股票"
$q1.Cells.Item(2, 3).Style = "Normal"
$q1.Cells.Item(2, 4).Value = "'12.82"
$q1.Cells.Item(2, 4).Style = "Normal"
$q1.Cells.Item(2, 5).Value = "'80.93"
$q1.Cells.Item(2, 5).Style = "Normal"
$q1.Cells.Item(2, 6).Value = "'4.76"
$q1.Cells.Item(2, 6).Style = "Normal"
$q1.Cells.Item(2, 7).Value = "'0.6102"
$q1.Cells.Item(2, 7).Style = "Normal"
$q1.Cells.Item(2, 8).Value = 6

# row 3
$q1.Cells.Item(3, 1).Value = 1
$q1.Cells.Item(3, 2).Value = "'004263"
$q1.Cells.Item(3, 2).Style = "Normal"
$q1.Cells.Item(3, 3).Value = "'华安沪港深机会灵活配置混合"
$q1.Cells.Item(3, 3).Style = "Normal"
$q1.Cells.Item(3, 4).Value = "'13.52"
$q1.Cells.Item(3, 4).Style = "Normal"
$q1.Cells.Item(3, 5).Value = "'94.79"
$q1.Cells.Item(3, 5).Style = "Normal"
$q1.Cells.Item(3, 6).Value = "'3.95"
$q1.Cells.Item(3, 6).Style = "Normal"
$q1.Cells.Item(3, 7).Value = "'0.5340"
$q1.Cells.Item(3, 7).Style = "Normal"
$q1.Cells.Item(3, 8).Value = 9

# row 4
$q1.Cells.Item(4, 1).Value = 2
$q1.Cells.Item(4, 2).Value = "'040011"
$q1.Cells.Item(4, 2).Style = "Normal"
$q1.Cells.Item(4, 3).Value = "'华安核心混合"
$q1.Cells.Item(4, 3).Style = "Normal"
$q1.Cells.Item(4, 4).Value = "'9.22"
$q1.Cells.Item(4, 4).Style = "Normal"
$q1.Cells.Item(4, 5).Value = "'88.18"
$q1.Cells.Item(4, 5).Style = "Normal"
$q1.Cells.Item(4, 6).Value = "'4.73"
$q1.Cells.Item(4, 6).Style = "Normal"
$q1.Cells.Item(4, 7).Value = "'0.4361"
$q1.Cells.Item(4, 7).Style = "Normal"
$q1.Cells.Item(4, 8).Value = 7

# row 5
$q1.Cells.Item(5, 1).Value = 3
$q1.Cells.Item(5, 2).Value = "'233009"
$q1.Cells.Item(5, 2).Style = "Normal"
$q1.Cells.Item(5, 3).Value = "'大摩多因子精选策略混合"
$q1.Cells.Item(5, 3).Style = "Normal"
$q1.Cells.Item(5, 4).Value = "'6.77"
$q1.Cells.Item(5, 4).Style = "Normal"
$q1.Cells.Item(5, 5).Value = "'89.73"
$q1.Cells.Item(5, 5).Style = "Normal"
$q1.Cells.Item(5, 6).Value = "'1.21"
$q1.Cells.Item(5, 6).Style = "Normal"
$q1.Cells.Item(5, 7).Value = "'0.0819"
$q1.Cells.Item(5, 7).Style = "Normal"
$q1.Cells.Item(5, 8).Value = 4

# row 6
$q1.Cells.Item(6, 1).Value = 4
$q1.Cells.Item(6, 2).Value = "'011231"
$q1.Cells.Item(6, 2).Style = "Normal"
$q1.Cells.Item(6, 3).Value = "'光大保德信锦弘混合A"
$q1.Cells.Item(6, 3).Style = "Normal"
$q1.Cells.Item(6, 4).Value = "'4.13"
$q1.Cells.Item(6, 4).Style = "Normal"
$q1.Cells.Item(6, 5).Value = "'20.96"
$q1.Cells.Item(6, 5).Style = "Normal"
$q1.Cells.Item(6, 6).Value = "'0.72"
$q1.Cells.Item(6, 6).Style = "Normal"
$q1.Cells.Item(6, 7).Value = "'0.0297"
$q1.Cells.Item(6, 7).Style = "Normal"
$q1.Cells.Item(6, 8).Value = 6

# row 7
$q1.Cells.Item(7, 1).Value = 5
$q1.Cells.Item(7, 2).Value = "'011232"
$q1.Cells.Item(7, 2).Style = "Normal"
$q1.Cells.Item(7, 3).Value = "'光大保德信锦弘混合C"
$q1.Cells.Item(7, 3).Style = "Normal"
$q1.Cells.Item(7, 4).Value = "'1.29"
$q1.Cells.Item(7, 4).Style = "Normal"
$q1.Cells.Item(7, 5).Value = "'20.96"
$q1.Cells.Item(7, 5).Style = "Normal"
$q1.Cells.Item(7, 6).Value = "'0.72"
$q1.Cells.Item(7, 6).Style = "Normal"
$q1.Cells.Item(7, 7).Value = "'0.0093"
$q1.Cells.Item(7, 7).Style = "Normal"
$q1.Cells.Item(7, 8).Value = 6

# row 8
$q1.Cells.Item(8, 1).Value = 6
$q1.Cells.Item(8, 2).Value = "'005021"
$q1.Cells.Item(8, 2).Style = "Normal"
$q1.Cells.Item(8, 3).Value = "'渤海汇金量化汇盈灵活配置混合"
$q1.Cells.Item(8, 3).Style = "Normal"
$q1.Cells.Item(8, 4).Value = "'0.02"
$q1.Cells.Item(8, 4).Style = "Normal"
$q1.Cells.Item(8, 5).Value = "'92.66"
$q1.Cells.Item(8, 5).Style = "Normal"
$q1.Cells.Item(8, 6).Value = "'1.34"
$q1.Cells.Item(8, 6).Style = "Normal"
$q1.Cells.Item(8, 7).Value = "'0.0003"
$q1.Cells.Item(8, 7).Style = "Normal"
$q1.Cells.Item(8, 8).Value = 3

# --- Step 3: update the "总计" (Total) summary sheet: insert a new top data row for 2022-Q1 ---
$total = $wb.Worksheets.Item("总计")

# Shift existing data rows (2021-Q4 / 2021-Q3) down by one, opening up row 2
$total.Rows.Item(2).Insert()

# The insert carries over neighbouring (header-row) formatting onto the blank row - strip it
$total.Range("A2:D2").ClearFormats()

# Re-apply the index-column style (bordered/centered) used throughout column A, copied from a sibling cell
$total.Cells.Item(3, 1).Copy($total.Cells.Item(2, 1))

# Fill in the new 2022-Q1 summary row
$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q1"
$total.Cells.Item(2, 3).Value = 7
$total.Cells.Item(2, 4).Value = 1.7

# Renumber the running index for the rows that shifted down
$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(4, 1).Value = 2

# Worksheets.Copy() leaves the newly inserted sheet as the active tab; restore the
# original active sheet (the diff doesn't touch 2021-Q3's sheetView/tabSelected).
$wb.Worksheets.Item("2021-Q3").Activate()

